$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A29").Value = "2025Q3"
$ws.Range("B29").Value = "2025Q3"
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 0.3442340791738382
